$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'womens tights tall',
    'womens tights tall long',
    'womens tights with pockets',
    'womens tights xl',
    'womens tommie copper',
    'womens trail pants',
    'womens trail runners',
    'womens travel pants quick dry',
    'womens tree leggings',
    'womens triathlon shorts',
    'womens under armour compression pants',
    'womens under armour workout',
    'womens underarmour basketball shorts',
    'womens underwear days of the week',
    'womens underwear fitness',
    'womens warm pants',
    'womens warm running leggings',
    'womens warm tights',
    'womens warm tights and leggings',
    'womens warmup pants',
    'womens weather pants',
    'womens white running tights',
    'womens wind pant',
    'womens winter running',
    'womens winter running skirt',
    'womens winter tennis skirts',
    'womens woekout pants',
    'womens working pants',
    'womens workout capri',
    'womens workout capri leggings',
    'womens workout capri pants',
    'womens workout capris',
    'womens workout capris high waist',
    'womens workout clothes shorts',
    'womens workout clothes winter',
    'womens workout clothing',
    'womens workout knee brace',
    'womens workout legging capris',
    'womens workout leggings',
    'womens workout leggings capri',
    'womens workout leggings high waist',
    'womens workout leggings lululemon',
    'womens workout leggings with design',
    'womens workout pants capri',
    'womens workout pants lululemon',
    'womens workout tights',
    'womens workout tights and leggings',
    'womens workout winter jacket',
    'womens xs leggings',
    'womens xs tights',
    'womens yoga capris',
    'womens yoga leggings capri',
    'womens yoga pants capri',
    'womens yoga pants high waist',
    'womens yoga pants long',
    'womens yoga pants long length',
    'womes running pants',
    'wommens clothes',
    'wonen tights',
    'wonens basketball',
    'wonens compression leggings',
    'wonens workout leggings',
    'wooden back support',
    'wooden bicycle stand',
    'wooden bike stand',
    'wooden man martial arts',
    'wooden zebra',
    'work capris for women',
    'work out capris',
    'work out clothing',
    'work rings for women',
    'work tights for women',
    'working leggings',
    'working leggings for women',
    'working support back brace',
    'workout abdominal support',
    'workout apparel women',
    'workout ball 65',
    'workout bottoms for women',
    'workout capri',
    'workout capri for women',
    'workout capri leggings',
    'workout capri leggings for women',
    'workout capri leggings for women high waist',
    'workout capris',
    'workout capris for women',
    'workout capris high waist',
    'workout capris leggings',
    'workout capris with pockets for women',
    'workout chart for home gym',
    'workout clothes capris',
    'workout clothes for cold weather',
    'workout clothes for women high waist',
    'workout clothes for women leggings',
    'workout clothes for women modest',
    'workout clothes for women winter',
    'workout clothes leggings',
    'workout clothing',
    'workout clothing for women',
    'workout clothing women'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
